# Add data for 2022-06-03: the "through May 25" snapshot rolls forward to
# "through May 26" (sheet name + header label), and the per-neighborhood /
# per-month counts pick up the carjackings that happened on 2022-06-03 in
# each year's running "through May 26" window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the matching header label in row 1.
$ws.Name = "Through 2022-05-26"
$ws.Range("B1").Value = "May 2022 (through May 26)"

# Updated / newly added counts (row = neighborhood, column = month).
$ws.Range("G3").Value  = 9    # Austin, May 2021
$ws.Range("Q3").Value  = 5    # Austin, May 2019
$ws.Range("V3").Value  = 3    # Austin, May 2018

$ws.Range("B5").Value  = 3    # Garfield Park, May 2022 (through May 26)
$ws.Range("Q5").Value  = 3    # Garfield Park, May 2019
$ws.Range("AF5").Value = 9    # Garfield Park, May 2016

$ws.Range("V23").Value  = 1   # Grand Crossing, May 2018 (new)
$ws.Range("AK23").Value = 1   # Grand Crossing, May 2015 (new)

$ws.Range("B24").Value = 4    # Grand Boulevard, May 2022 (through May 26)
$ws.Range("L24").Value = 2    # Grand Boulevard, May 2020

$ws.Range("G28").Value  = 3   # West Town, May 2021
$ws.Range("AA28").Value = 1   # West Town, May 2017 (new)

$ws.Range("G55").Value = 2    # Bucktown, May 2021 (new)

$ws.Range("AA80").Value = 1   # Oakland, May 2017 (new)

$ws.Range("AA91").Value = 1   # Washington Park, May 2017 (new)
